# Applies updated PSSM score values (B2:K21) per commit "updates with supplemental figures".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'double[,]' 20,10
$data[0,0] = -19.4817659333861
$data[0,1] = 2.348506702695107
$data[0,2] = -19.4817659333861
$data[0,3] = -19.4817659333861
$data[0,4] = -19.4817659333861
$data[0,5] = -19.4817659333861
$data[0,6] = -19.4817659333861
$data[0,7] = -19.4817659333861
$data[0,8] = -19.4817659333861
$data[0,9] = -19.4817659333861
$data[1,0] = -19.4817659333861
$data[1,1] = -19.4817659333861
$data[1,2] = -19.4817659333861
$data[1,3] = -19.4817659333861
$data[1,4] = -19.4817659333861
$data[1,5] = -19.4817659333861
$data[1,6] = -19.4817659333861
$data[1,7] = 4.321926222901729
$data[1,8] = -19.4817659333861
$data[1,9] = -19.4817659333861
$data[2,0] = -19.4817659333861
$data[2,1] = 2.117649811039565
$data[2,2] = 2.786691828509597
$data[2,3] = -19.4817659333861
$data[2,4] = 2.55935075022294
$data[2,5] = -19.4817659333861
$data[2,6] = 1.768796770253058
$data[2,7] = -19.4817659333861
$data[2,8] = 2.316418287925416
$data[2,9] = -19.4817659333861
$data[3,0] = -19.4817659333861
$data[3,1] = 1.005600143800447
$data[3,2] = -19.4817659333861
$data[3,3] = -19.4817659333861
$data[3,4] = -19.4817659333861
$data[3,5] = 1.924607774994609
$data[3,6] = -19.4817659333861
$data[3,7] = -19.4817659333861
$data[3,8] = -19.4817659333861
$data[3,9] = -19.4817659333861
$data[4,0] = -19.4817659333861
$data[4,1] = -19.4817659333861
$data[4,2] = -19.4817659333861
$data[4,3] = -19.4817659333861
$data[4,4] = -19.4817659333861
$data[4,5] = -19.4817659333861
$data[4,6] = -19.4817659333861
$data[4,7] = -19.4817659333861
$data[4,8] = -19.4817659333861
$data[4,9] = -19.4817659333861
$data[5,0] = 2.961630165611308
$data[5,1] = -19.4817659333861
$data[5,2] = -19.4817659333861
$data[5,3] = -19.4817659333861
$data[5,4] = -19.4817659333861
$data[5,5] = -19.4817659333861
$data[5,6] = -19.4817659333861
$data[5,7] = -19.4817659333861
$data[5,8] = -19.4817659333861
$data[5,9] = -19.4817659333861
$data[6,0] = -19.4817659333861
$data[6,1] = -19.4817659333861
$data[6,2] = -19.4817659333861
$data[6,3] = 2.808109148668288
$data[6,4] = -19.4817659333861
$data[6,5] = -19.4817659333861
$data[6,6] = -19.4817659333861
$data[6,7] = -19.4817659333861
$data[6,8] = -19.4817659333861
$data[6,9] = -19.4817659333861
$data[7,0] = 3.609984028031273
$data[7,1] = -19.4817659333861
$data[7,2] = -19.4817659333861
$data[7,3] = -19.4817659333861
$data[7,4] = -19.4817659333861
$data[7,5] = -19.4817659333861
$data[7,6] = -19.4817659333861
$data[7,7] = -19.4817659333861
$data[7,8] = -19.4817659333861
$data[7,9] = -19.4817659333861
$data[8,0] = -19.4817659333861
$data[8,1] = -19.4817659333861
$data[8,2] = -19.4817659333861
$data[8,3] = -19.4817659333861
$data[8,4] = -19.4817659333861
$data[8,5] = -19.4817659333861
$data[8,6] = -19.4817659333861
$data[8,7] = -19.4817659333861
$data[8,8] = -19.4817659333861
$data[8,9] = 2.311059701656554
$data[9,0] = -19.4817659333861
$data[9,1] = -19.4817659333861
$data[9,2] = -19.4817659333861
$data[9,3] = 2.072311827065299
$data[9,4] = -19.4817659333861
$data[9,5] = 2.54378899995148
$data[9,6] = -19.4817659333861
$data[9,7] = -19.4817659333861
$data[9,8] = -19.4817659333861
$data[9,9] = 1.296100169364725
$data[10,0] = -19.4817659333861
$data[10,1] = -19.4817659333861
$data[10,2] = -19.4817659333861
$data[10,3] = -19.4817659333861
$data[10,4] = -19.4817659333861
$data[10,5] = -19.4817659333861
$data[10,6] = -19.4817659333861
$data[10,7] = -19.4817659333861
$data[10,8] = -19.4817659333861
$data[10,9] = -19.4817659333861
$data[11,0] = -19.4817659333861
$data[11,1] = -19.4817659333861
$data[11,2] = -19.4817659333861
$data[11,3] = 1.790228746865409
$data[11,4] = -19.4817659333861
$data[11,5] = -19.4817659333861
$data[11,6] = -19.4817659333861
$data[11,7] = -19.4817659333861
$data[11,8] = 2.188795666882015
$data[11,9] = 1.609070779931526
$data[12,0] = -19.4817659333861
$data[12,1] = -19.4817659333861
$data[12,2] = 1.647341591616368
$data[12,3] = -19.4817659333861
$data[12,4] = -19.4817659333861
$data[12,5] = -19.4817659333861
$data[12,6] = -19.4817659333861
$data[12,7] = -19.4817659333861
$data[12,8] = -19.4817659333861
$data[12,9] = 2.141774902701559
$data[13,0] = -19.4817659333861
$data[13,1] = -19.4817659333861
$data[13,2] = -0.2301147904235224
$data[13,3] = -19.4817659333861
$data[13,4] = -19.4817659333861
$data[13,5] = -19.4817659333861
$data[13,6] = -19.4817659333861
$data[13,7] = -19.4817659333861
$data[13,8] = -19.4817659333861
$data[13,9] = -19.4817659333861
$data[14,0] = -19.4817659333861
$data[14,1] = -19.4817659333861
$data[14,2] = -19.4817659333861
$data[14,3] = -19.4817659333861
$data[14,4] = -19.4817659333861
$data[14,5] = -19.4817659333861
$data[14,6] = -19.4817659333861
$data[14,7] = -19.4817659333861
$data[14,8] = 2.257075528145628
$data[14,9] = -19.4817659333861
$data[15,0] = -19.4817659333861
$data[15,1] = 0.8063758094609516
$data[15,2] = 0.113231058990873
$data[15,3] = -19.4817659333861
$data[15,4] = -19.4817659333861
$data[15,5] = -19.4817659333861
$data[15,6] = 0.2126888802919562
$data[15,7] = -19.4817659333861
$data[15,8] = 1.253675226961352
$data[15,9] = -19.4817659333861
$data[16,0] = -19.4817659333861
$data[16,1] = -19.4817659333861
$data[16,2] = -19.4817659333861
$data[16,3] = -19.4817659333861
$data[16,4] = -19.4817659333861
$data[16,5] = -19.4817659333861
$data[16,6] = 0.4603261319065275
$data[16,7] = -19.4817659333861
$data[16,8] = 1.720320216116217
$data[16,9] = -19.4817659333861
$data[17,0] = -19.4817659333861
$data[17,1] = -19.4817659333861
$data[17,2] = 1.818526851571873
$data[17,3] = -19.4817659333861
$data[17,4] = -19.4817659333861
$data[17,5] = -19.4817659333861
$data[17,6] = 1.771581259876474
$data[17,7] = -19.4817659333861
$data[17,8] = -19.4817659333861
$data[17,9] = -19.4817659333861
$data[18,0] = -19.4817659333861
$data[18,1] = 1.775442910349884
$data[18,2] = 2.171700651377595
$data[18,3] = -19.4817659333861
$data[18,4] = 3.818191197597642
$data[18,5] = -19.4817659333861
$data[18,6] = 2.163419716190238
$data[18,7] = -19.4817659333861
$data[18,8] = -19.4817659333861
$data[18,9] = 2.355654501986832
$data[19,0] = -19.4817659333861
$data[19,1] = 1.760061893553971
$data[19,2] = -19.4817659333861
$data[19,3] = 2.414682147949911
$data[19,4] = -19.4817659333861
$data[19,5] = 3.374685622083546
$data[19,6] = 2.623761858368571
$data[19,7] = -19.4817659333861
$data[19,8] = -19.4817659333861
$data[19,9] = -19.4817659333861

$ws.Range("B2:K21").Value = $data
Write-Output "Updated B2:K21 with new PSSM values"
